$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update employee ids (col A) to the new numbering scheme ---
$newIds = @(81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98)
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}

# --- Every weekoff (col B) now carries the same corrected date ---
$ws.Range("B2:B19").Value = "6/2/2024"

# --- The id column no longer needs its old Arial/bordered look, drop back to default ---
$ws.Range("A2:A19").ClearFormats()

# --- Leave the cursor where the next rows would be entered ---
$ws.Range("A20:A21").Select() | Out-Null
